$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 - Subject 6
$ws.Range("B8").Value = "Female"
$ws.Range("C8").Value = 22
$ws.Range("D8").Value = "Psychology"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "Mobile; Game Console"
$ws.Range("G8").Value = "Joystick"
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = "G"

# Row 9 - Subject 7
$ws.Range("B9").Value = "Male"
$ws.Range("C9").Value = 21
$ws.Range("D9").Value = "CS"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = "Tablet"
$ws.Range("G9").Value = "Other"
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 2
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = "F"

# Row 10 - Subject 8
$ws.Range("B10").Value = "Male"
$ws.Range("C10").Value = 20
$ws.Range("D10").Value = "IMGD/CS"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = "Desktop: Game Console"
$ws.Range("G10").Value = "Keyboard/Mouse"
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 4
$ws.Range("K10").Value = "G"

# Row 11 - Subject 9
$ws.Range("B11").Value = "Male"
$ws.Range("C11").Value = 20
$ws.Range("D11").Value = "Mechanical Eng."
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = "DeskTop; Mobile; Game Console"
$ws.Range("G11").Value = "Keyboard/Mouse; Joystick; "
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 4
$ws.Range("K11").Value = "F"

# Rows with long wrapped text grow taller (same rows the diff marks ht="30")
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30

# Update the active selection to match the final state of the edit
$ws.Range("J11").Select()
